$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 17999.777
$ws.Range("I21").Value = 17999.777
$ws.Range("K21").Value = 17999.777
$ws.Range("M21").Value = -17531.777
$ws.Range("H23").Value = 17999.777
$ws.Range("I23").Value = 17999.777
$ws.Range("K23").Value = 17999.777
$ws.Range("M23").Value = -17765.777
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("H80").Value = 1354.8182
$ws.Range("J80").Value = 1431.625
$ws.Range("L80").Value = 4294.875
$ws.Range("N80").Value = -6290.875
$ws.Range("H83").Value = 1354.8182
$ws.Range("J83").Value = 1431.625
$ws.Range("L83").Value = 12884.625
$ws.Range("N83").Value = -22868.625
$ws.Range("H88").Value = 701.4
$ws.Range("J88").Value = 782.5
$ws.Range("L88").Value = 782.5
$ws.Range("N88").Value = -1594.5
$ws.Range("H91").Value = 701.4
$ws.Range("J91").Value = 782.5
$ws.Range("L91").Value = 782.5
$ws.Range("N91").Value = -3590.5
$ws.Range("H100").Value = 7899.2
$ws.Range("I100").Value = 7748
$ws.Range("K100").Value = 7748
$ws.Range("M100").Value = -7207
$ws.Range("H112").Value = 3208
$ws.Range("J112").Value = 3208
$ws.Range("L112").Value = 9624
$ws.Range("N112").Value = -11840
$ws.Range("H127").Value = 4922.5713
$ws.Range("I127").Value = 4909.6665
$ws.Range("K127").Value = 14728.9995
$ws.Range("M127").Value = -9768.999500000002
$ws.Range("H137").Value = 2739.8948
$ws.Range("I137").Value = 1545.9166
$ws.Range("K137").Value = 4637.7498
$ws.Range("M137").Value = -2087.7498
$ws.Range("H138").Value = 1822.6666
$ws.Range("I138").Value = 1583.4
$ws.Range("J138").Value = 2121.75
$ws.Range("K138").Value = 4750.200000000001
$ws.Range("L138").Value = 6365.25
$ws.Range("M138").Value = 389.7999999999993
$ws.Range("N138").Value = -16645.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5568.25
$ws.Range("I2").Value = 2059.6
$ws.Range("K2").Value = 2059.6
$ws.Range("M2").Value = -1946.6
$ws.Range("H61").Value = 5111.25
$ws.Range("I61").Value = 4412.857
$ws.Range("K61").Value = 4412.857
$ws.Range("M61").Value = -4200.857
$ws.Range("H74").Value = 850.5625
$ws.Range("I74").Value = 867.26666
$ws.Range("J74").Value = 600
$ws.Range("K74").Value = 867.26666
$ws.Range("L74").Value = 600
$ws.Range("M74").Value = 6.733339999999998
$ws.Range("N74").Value = -2348
$ws.Range("H77").Value = 850.5625
$ws.Range("I77").Value = 867.26666
$ws.Range("K77").Value = 4336.3333
$ws.Range("L77").Value = 3000
$ws.Range("M77").Value = 31.66669999999976
$ws.Range("N77").Value = -11736
$ws.Range("H95").Value = 24208
$ws.Range("J95").Value = 24208
$ws.Range("L95").Value = 24208
$ws.Range("N95").Value = -29700
$ws.Range("H116").Value = 5568.25
$ws.Range("I116").Value = 2059.6
$ws.Range("K116").Value = 2059.6
$ws.Range("M116").Value = 234.4000000000001
$ws.Range("H132").Value = 3924.65
$ws.Range("I132").Value = 3249.611
$ws.Range("K132").Value = 9748.832999999999
$ws.Range("M132").Value = -7218.832999999999
$ws.Range("H136").Value = 5111.25
$ws.Range("I136").Value = 4412.857
$ws.Range("K136").Value = 13238.571
$ws.Range("M136").Value = -10688.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5568.25
$ws.Range("I3").Value = 2059.6
$ws.Range("K3").Value = 2059.6
$ws.Range("M3").Value = -1945.6
$ws.Range("H80").Value = 215.21053
$ws.Range("I80").Value = 184.5
$ws.Range("J80").Value = 223.4
$ws.Range("K80").Value = 184.5
$ws.Range("L80").Value = 223.4
$ws.Range("M80").Value = 813.5
$ws.Range("N80").Value = -2219.4
$ws.Range("H83").Value = 215.21053
$ws.Range("I83").Value = 184.5
$ws.Range("J83").Value = 223.4
$ws.Range("K83").Value = 922.5
$ws.Range("L83").Value = 1117
$ws.Range("M83").Value = 4069.5
$ws.Range("N83").Value = -11101

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6100.6387
$ws.Range("I31").Value = 2344.0588
$ws.Range("K31").Value = 2344.0588
$ws.Range("M31").Value = -2049.0588
$ws.Range("H34").Value = 6100.6387
$ws.Range("I34").Value = 2344.0588
$ws.Range("K34").Value = 2344.0588
$ws.Range("M34").Value = -2142.0588
$ws.Range("H52").Value = 190000
$ws.Range("J52").Value = 250000
$ws.Range("L52").Value = 250000
$ws.Range("N52").Value = -250588
$ws.Range("H99").Value = 5110.091
$ws.Range("I99").Value = 4606.8887
$ws.Range("K99").Value = 4606.8887
$ws.Range("M99").Value = -3108.8887
$ws.Range("H107").Value = 621.38464
$ws.Range("I107").Value = 594
$ws.Range("K107").Value = 594
$ws.Range("M107").Value = 1326
$ws.Range("H126").Value = 5110.091
$ws.Range("I126").Value = 4606.8887
$ws.Range("K126").Value = 13820.6661
$ws.Range("M126").Value = -11350.6661
$ws.Range("H132").Value = 2442.15
$ws.Range("I132").Value = 1908.2222
$ws.Range("J132").Value = 7247.5
$ws.Range("K132").Value = 5724.6666
$ws.Range("L132").Value = 21742.5
$ws.Range("M132").Value = -3194.6666
$ws.Range("N132").Value = -26802.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1918.8
$ws.Range("I5").Value = 863
$ws.Range("K5").Value = 2589
$ws.Range("M5").Value = -2477
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("H107").Value = 390.66666
$ws.Range("J107").Value = 364.66666
$ws.Range("L107").Value = 1093.99998
$ws.Range("N107").Value = -4933.999980000001
$ws.Range("H135").Value = 1918.8
$ws.Range("I135").Value = 863
$ws.Range("K135").Value = 7767
$ws.Range("M135").Value = -5232
$ws.Range("H137").Value = 2663.1428
$ws.Range("J137").Value = 2723.5
$ws.Range("L137").Value = 8170.5
$ws.Range("N137").Value = -18370.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7661.5386
$ws.Range("I113").Value = 6085.7144
$ws.Range("J113").Value = 9500
$ws.Range("K113").Value = 6085.7144
$ws.Range("L113").Value = 9500
$ws.Range("M113").Value = -3915.7144
$ws.Range("N113").Value = -13840
$ws.Range("H122").Value = 2376.25
$ws.Range("I122").Value = 1604.7142
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 4814.142599999999
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -2364.142599999999
$ws.Range("N122").Value = -28231
$ws.Range("H132").Value = 142563.25
$ws.Range("I132").Value = 220604.6
$ws.Range("K132").Value = 661813.8
$ws.Range("M132").Value = -659283.8
$ws.Range("H135").Value = 232500
$ws.Range("J135").Value = 232500
$ws.Range("L135").Value = 232500
$ws.Range("N135").Value = -242640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7050.3887
$ws.Range("I40").Value = 6995.0835
$ws.Range("K40").Value = 6995.0835
$ws.Range("M40").Value = -6859.0835
$ws.Range("H55").Value = 851.7692
$ws.Range("I55").Value = 970.63635
$ws.Range("J55").Value = 764.6
$ws.Range("K55").Value = 970.63635
$ws.Range("L55").Value = 764.6
$ws.Range("M55").Value = -797.63635
$ws.Range("N55").Value = -1110.6
$ws.Range("H68").Value = 8916.5
$ws.Range("J68").Value = 8916.5
$ws.Range("L68").Value = 8916.5
$ws.Range("N68").Value = -10414.5
$ws.Range("H71").Value = 8916.5
$ws.Range("J71").Value = 8916.5
$ws.Range("L71").Value = 44582.5
$ws.Range("N71").Value = -52070.5
$ws.Range("H122").Value = 3995.6667
$ws.Range("I122").Value = 3995.6667
$ws.Range("K122").Value = 11987.0001
$ws.Range("M122").Value = -9537.000100000001
$ws.Range("H135").Value = 89000
$ws.Range("J135").Value = 89000
$ws.Range("L135").Value = 89000
$ws.Range("N135").Value = -99140
$ws.Range("H137").Value = 419500
$ws.Range("I137").Value = 89000
$ws.Range("K137").Value = 89000
$ws.Range("M137").Value = -83900
